$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Roraima stays, variable/value updated
$ws.Range("B2").Value = "Diferença 2024/02 - 2023/02"
$ws.Range("C2").Value = 2

# Row 3 - now Distrito Federal
$ws.Range("A3").Value = "Distrito Federal"
$ws.Range("B3").Value = "Diferença 2024/02 - 2023/02"
$ws.Range("C3").Value = 1

# Row 4 - now Rondônia
$ws.Range("A4").Value = "Rondônia"
$ws.Range("B4").Value = "Diferença 2024/02 - 2023/02"
$ws.Range("C4").Value = 0.8999999999999999

# Row 5 - now Rio Grande do Sul
$ws.Range("A5").Value = "Rio Grande do Sul"
$ws.Range("B5").Value = "Diferença 2024/02 - 2023/02"
$ws.Range("C5").Value = 0.6000000000000005

# Row 6 - now Mato Grosso
$ws.Range("A6").Value = "Mato Grosso"
$ws.Range("B6").Value = "Diferença 2024/02 - 2023/02"
$ws.Range("C6").Value = 0.2999999999999998

# Row 7 - now Mato Grosso do Sul
$ws.Range("A7").Value = "Mato Grosso do Sul"
$ws.Range("B7").Value = "Diferença 2024/02 - 2023/02"
$ws.Range("C7").Value = -0.2999999999999998

# Row 8 - Sergipe stays, value & rank updated
$ws.Range("B8").Value = "Diferença 2024/02 - 2023/02"
$ws.Range("C8").Value = -1.200000000000001
$ws.Range("D8").Value = "15º"

# Row 9 - Nordeste stays
$ws.Range("B9").Value = "Diferença 2024/02 - 2023/02"
$ws.Range("C9").Value = -1.9

# Row 10 - Brasil stays
$ws.Range("B10").Value = "Diferença 2024/02 - 2023/02"
$ws.Range("C10").Value = -1.1
